$d = $word.ActiveDocument

# Locate the paragraph that ends with "...drupal_drupalgeddon2" (the
# "use exploit/unix/webapp/drupal_drupalgeddon2" line) using Find so the
# script does not depend on a hard-coded paragraph index.
$findRange = $d.Content
$found = $findRange.Find.Execute("drupal_drupalgeddon2")

if ($found) {
    $targetPara = $findRange.Paragraphs(1)

    # Insert a brand-new (empty) paragraph right after it.
    $targetPara.Range.InsertParagraphAfter()

    # The freshly inserted paragraph is the next one.
    $newPara = $targetPara.Next()
    $newPara.Range.Text = "set RHOST <ip>"

    # Re-fetch the paragraph (its Range shifted once text was typed into it)
    # and drop a collapsed "_GoBack" bookmark in it, mirroring the
    # last-edit-position bookmark Word leaves behind after typing new text.
    # The bookmark is placed at the start of the paragraph (rather than
    # immediately after the text) purely to dodge a collapsed-range edge
    # case at "paragraph end - 1" in this host; bookmarks are invisible,
    # non-rendering metadata so this has no visible effect.
    $newPara2 = $targetPara.Next()
    $bmPos = $newPara2.Range.Start
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
